# Widen the AMAM table columns / overall table width to prevent wrapping.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Overall preferred table width: w:tblW w:w goes from 4933 -> 4946 (pct*50),
# and PreferredWidth (points-like units) = w:w / 20.
$t.PreferredWidth = 247.3

# Column (gridCol) widths, in points; OOXML twips = points * 20.
#   531  -> 428   =>  21.40
#   1647 -> 1327  =>  66.35
#   2870 -> 3039  => 151.95
#   2764 -> 3039  => 151.95
$t.Columns.Item(1).Width = 21.4
$t.Columns.Item(2).Width = 66.35
$t.Columns.Item(3).Width = 151.95
$t.Columns.Item(4).Width = 151.95
